$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-15"

# Update the December label cell (A13) to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-15)"

# Update December row (row 13) with new cumulative counts
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 46
$ws.Range("D13").Value = 52
$ws.Range("E13").Value = 34
$ws.Range("F13").Value = 25
$ws.Range("G13").Value = 77
$ws.Range("H13").Value = 114
$ws.Range("I13").Value = 64

# Update Total row (row 14) with new cumulative counts
$ws.Range("B14").Value = 307
$ws.Range("C14").Value = 609
$ws.Range("D14").Value = 873
$ws.Range("E14").Value = 716
$ws.Range("F14").Value = 559
$ws.Range("G14").Value = 1341
$ws.Range("H14").Value = 1757
$ws.Range("I14").Value = 1580
